# 디비 컬럼명 변경 및 ticket, reser, airline 디비 선언
# Update the "model" sheet (DB structure) per the commit:
#  - Ticket table: rename/expand columns, split price into 4 seat-class prices
#  - Reservation table: rename go/come ticket id + date columns
# sheet1 ("웹구조") is not semantically changed by the diff (only shared-string
# index shuffling caused by sharedStrings.xml edits elsewhere), so it is left
# untouched here.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("model")

# --- Ticket table header row (row 5): K is now "날짜", new L/M/N/O columns ---
$ws.Range("K5").Value = "날짜"
$ws.Range("L5").Value = "일등석"
$ws.Range("M5").Value = "비즈니스석"
$ws.Range("N5").Value = "프리미엄석"
$ws.Range("O5").Value = "일반석"

# --- Ticket table column-name row (row 6): renamed DB field names ---
$ws.Range("D6").Value = "airline_id"
$ws.Range("E6").Value = "departure_place"
$ws.Range("F6").Value = "arrival_place"
$ws.Range("G6").Value = "departure_airport"
$ws.Range("H6").Value = "arrival_airport"
$ws.Range("I6").Value = "departure_time"
$ws.Range("J6").Value = "arrival_time"
$ws.Range("K6").Value = "departure_date"
$ws.Range("L6").Value = "first_class_price"
$ws.Range("M6").Value = "business_class_price"
$ws.Range("N6").Value = "premium_price"
$ws.Range("O6").Value = "economy_price"

# --- Reservation table column-name row (row 10): renamed DB field names ---
$ws.Range("E10").Value = "go_ticket_id"
$ws.Range("F10").Value = "come_ticket_id"
$ws.Range("G10").Value = "reservation_date"

# --- Column widths for the newly populated / widened columns (best effort
#     match to the post-edit bestFit widths) ---
$ws.Columns.Item(7).ColumnWidth = 15.142857142857142
$ws.Columns.Item(8).ColumnWidth = 11.714285714285714
$ws.Columns.Item(9).ColumnWidth = 13.285714285714286
$ws.Columns.Item(10).ColumnWidth = 9.857142857142858
$ws.Columns.Item(11).ColumnWidth = 13.285714285714286
$ws.Columns.Item(12).ColumnWidth = 13.285714285714286
$ws.Columns.Item(13).ColumnWidth = 17.714285714285715
$ws.Columns.Item(14).ColumnWidth = 13.142857142857142
$ws.Columns.Item(15).ColumnWidth = 13.428571428571429

# --- View: selection moves to F7 (and any topLeftCell scroll pin is cleared) ---
$ws.Activate() | Out-Null
$ws.Range("F7").Select() | Out-Null
